# Applies the "Criação da tabela 1 e histogramas + comentários" edit:
#  1) Adds spell-check proofErr wrapping ("sq.ft", "feet") and splits
#     runs in the ZN (variable 2) description.
#  2) Splits the PTRATIO (variable 11) run into three runs and
#     capitalises "Proporção".
#  3) Adds spell-check proofErr wrapping ("Bk" x2) and grammar-check
#     proofErr wrapping (around "63)^") in the B (variable 12) equation
#     description.
#  4) Appends two empty paragraphs and a new paragraph introducing the
#     Figura 1 histograms after the LSTAT (variable 13) description.

$d = $word.ActiveDocument

function Wrap-Pkg($bodyInner) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParaIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) ZN paragraph (variable 2): "sq.ft" and "Square feet" proof marks.
# ---------------------------------------------------------------------
$znIdx = Find-ParaIndex $d "sq.ft"
$znXml = '<w:p>' +
    '<w:r><w:t>2</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> ZN: </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Proporção </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">de terreno residencial zoneada para lotes acima de 25.000 </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>sq.ft</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Square </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>fee</w:t></w:r>' +
    '<w:r><w:t>t</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> – Pés quadrados)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '</w:p>'
$d.Paragraphs.Item($znIdx).Range.InsertXML((Wrap-Pkg $znXml))

# ---------------------------------------------------------------------
# 2) PTRATIO paragraph (variable 11): split into three runs.
# ---------------------------------------------------------------------
$ptratioIdx = Find-ParaIndex $d "PTRATIO"
$ptratioXml = '<w:p>' +
    '<w:r><w:t>11</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> PTRATIO: </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Proporção </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">aluno-professor por bairro.  </w:t></w:r>' +
    '</w:p>'
$d.Paragraphs.Item($ptratioIdx).Range.InsertXML((Wrap-Pkg $ptratioXml))

# ---------------------------------------------------------------------
# 3) B paragraph (variable 12): "Bk" (x2) spell marks + grammar marks.
# ---------------------------------------------------------------------
$bIdx = Find-ParaIndex $d "B: O resultado"
$bXml = '<w:p>' +
    '<w:r><w:t>12</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> B: O resultado da equação $B=1000(</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Bk</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> - 0,</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>63</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '<w:r><w:t>^</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>2$ onde $</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Bk</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">$ é a proporção de negros por bairro.    </w:t></w:r>' +
    '</w:p>'
$d.Paragraphs.Item($bIdx).Range.InsertXML((Wrap-Pkg $bXml))

# ---------------------------------------------------------------------
# 4) Append two empty paragraphs + the new "Figura 1" paragraph after
#    the LSTAT paragraph (variable 13, the last paragraph of the body).
#
#    InsertXML-ing a range anchored at the very end of the document's
#    main story leaves the original trailing paragraph mark behind as
#    an extra empty paragraph after whatever was inserted, so a spare
#    placeholder paragraph is created first (InsertParagraphAfter) and
#    the genuine stray mark it leaves behind is trimmed afterwards.
# ---------------------------------------------------------------------
$lstatIdx = Find-ParaIndex $d "LSTAT"
$appendXml = '<w:p/><w:p/><w:p>' +
    '<w:r><w:t xml:space="preserve">A Figura 1 traz os Histogramas das variáveis a fim de tentar identificar o </w:t></w:r>' +
    '<w:r><w:t>comportamento</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> das </w:t></w:r>' +
    '<w:r><w:t>variáveis</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> e uma </w:t></w:r>' +
    '<w:r><w:t>possível</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> associação a algum modelo de </w:t></w:r>' +
    '<w:r><w:t>distribuição</w:t></w:r>' +
    '</w:p>'

$lstatRange = $d.Paragraphs.Item($lstatIdx).Range
$lstatRange.InsertParagraphAfter()
$placeholder = $d.Paragraphs.Item($lstatIdx + 1)
$placeholder.Range.InsertXML((Wrap-Pkg $appendXml))

# Trim the stray empty paragraph mark left dangling at the new story end.
$tailIdx = $d.Paragraphs.Count
$tail = $d.Paragraphs.Item($tailIdx).Range
$d.Range($tail.Start - 1, $tail.End).Delete()
